# Add team record columns (Wins/Losses/Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the existing header style (from AC1, style index 1:
# bold, bordered, centered) onto the three new header cells AD1:AF1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Data rows 2 through 53: every player row gets the same team record,
# 80 wins, 82 losses, 0 ties.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 80
    $ws.Cells.Item($r, 31).Value2 = 82
    $ws.Cells.Item($r, 32).Value2 = 0
}
